# Update the "想去人数" (want-to-go count) figures in column F on the
# "展览" sheet and the aggregated "全部类型" sheet, matching the latest
# scrape output (gh-pages data refresh at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 7002
$ws1.Range("F5").Value  = 457
$ws1.Range("F6").Value  = 157
$ws1.Range("F7").Value  = 6919
$ws1.Range("F8").Value  = 75
$ws1.Range("F9").Value  = 203
$ws1.Range("F10").Value = 1287
$ws1.Range("F13").Value = 410
$ws1.Range("F14").Value = 151
$ws1.Range("F15").Value = 18
$ws1.Range("F16").Value = 416
$ws1.Range("F17").Value = 51
$ws1.Range("F18").Value = 43
$ws1.Range("F19").Value = 18
$ws1.Range("F20").Value = 5285
$ws1.Range("F23").Value = 702
$ws1.Range("F24").Value = 222
$ws1.Range("F25").Value = 247

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 7002
$ws4.Range("F4").Value  = 64
$ws4.Range("F7").Value  = 6919
$ws4.Range("F8").Value  = 75
$ws4.Range("F9").Value  = 203
$ws4.Range("F10").Value = 1287
$ws4.Range("F11").Value = 22
$ws4.Range("F12").Value = 109
$ws4.Range("F14").Value = 151
$ws4.Range("F15").Value = 18
$ws4.Range("F16").Value = 416
$ws4.Range("F17").Value = 51
$ws4.Range("F18").Value = 43
$ws4.Range("F19").Value = 18
$ws4.Range("F21").Value = 5285
$ws4.Range("F22").Value = 46
$ws4.Range("F25").Value = 702
$ws4.Range("F26").Value = 222
$ws4.Range("F27").Value = 247
